# eFinance_Process/Data/Config.xlsx update
#
# Commit message: "Updated the BillingStatement sequence Created the
# GetAttachmentsSequence"
#
# Concretely this adds:
#   - Constants sheet: three new named timeout constants (TimeoutShort,
#     TimeoutMedium, TimeoutLong) used by the new GetAttachments sequence.
#   - Settings sheet: a SupportedFileFormats setting used by the new
#     GetAttachments.xaml workflow.

$wb = $excel.ActiveWorkbook

# --- Constants sheet: rows 12-14 (TimeoutShort / TimeoutMedium / TimeoutLong) ---
$wsConstants = $wb.Worksheets.Item("Constants")

$wsConstants.Range("A12").Value = "TimeoutShort"
$wsConstants.Range("A13").Value = "TimeoutMedium"
$wsConstants.Range("A14").Value = "TimeoutLong"

$wsConstants.Range("B12").Value = 5000
$wsConstants.Range("B13").Value = 30000
$wsConstants.Range("B14").Value = 120000

$wsConstants.Range("C12").Value = "Timeout short value in milliseconds, for activities which are likely to fail. Must be integer"
$wsConstants.Range("C13").Value = "Timeout medium value in milliseconds. Must be integer"
$wsConstants.Range("C14").Value = "Timeout short value in milliseconds, for slow apps. Must be integer"

# Leave the Constants sheet's last selection on C14, matching the cell the
# author was last working in.
$wsConstants.Range("C14").Select() | Out-Null

# --- Settings sheet: row 7 (SupportedFileFormats) ---
$wsSettings = $wb.Worksheets.Item("Settings")

$wsSettings.Range("A7").Value = "SupportedFileFormats"
$wsSettings.Range("B7").Value = '{".pdf",".jpg",".jpeg"}'
$wsSettings.Range("C7").Value = "List of all supported file formats for the GetAttachments.xaml workflow"

# Settings is the active/visible sheet, leave the selection on the new row.
$wsSettings.Range("A7").Select() | Out-Null
